$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 58: becomes the "CPOF" (CAFÉ PERGAMINO ORGÁNICO FAIRTRADE) data row ---
# row 58 already carries s="3" customFormat="1" at the row level (unchanged); we
# just need to fill in the values and mark the B:G cells bold to reproduce the
# explicit per-cell s="3" seen in the target.
$ws.Range("B58").Value2 = 1
$ws.Range("C58").Value2 = 2007
$ws.Range("D58").Value2 = 1001
$ws.Range("E58").Value2 = "CPOF"
$ws.Range("F58").Value2 = "CAFÉ PERGAMINO ORGÁNICO FAIRTRADE"
$ws.Range("G58").Value2 = 1
$ws.Range("B58:G58").Font.Bold = $true

# --- Row 59: "CPF" (CAFÉ PERGAMINO FAIRTRADE) data row, plain/default style ---
$ws.Range("B59").Value2 = 2
$ws.Range("C59").Value2 = 2007
$ws.Range("D59").Value2 = 1002
$ws.Range("E59").Value2 = "CPF"
$ws.Range("F59").Value2 = "CAFÉ PERGAMINO FAIRTRADE"
$ws.Range("G59").Value2 = 1

# --- Row 60: "CO" (CAFÉ ORGÁNICO) data row, plain/default style ---
$ws.Range("B60").Value2 = 3
$ws.Range("C60").Value2 = 2007
$ws.Range("D60").Value2 = 1003
$ws.Range("E60").Value2 = "CO"
$ws.Range("F60").Value2 = "CAFÉ ORGÁNICO"
$ws.Range("G60").Value2 = 1

# --- Row 61: new "COSECHAS" header/class row (bold cells, plain row) ---
$ws.Range("B61").Value2 = 0
$ws.Range("C61").Value2 = 2008
$ws.Range("D61").Value2 = 10
$ws.Range("E61").Value2 = "COSECHAS"
$ws.Range("F61").Value2 = ".:::.COSECHAS.:::."
$ws.Range("G61").Value2 = 0
$ws.Range("B61:G61").Font.Bold = $true

# --- Row 63: the old "ROLES" row (B63:G63) moves down to row 70; clear it here ---
$ws.Range("B63:G63").ClearContents()

# --- Row 70: re-create the "ROLES" row that used to live at row 63 ---
$ws.Range("B70").Value2 = 0
$ws.Range("C70").Value2 = 2038
$ws.Range("D70").Value2 = 10
$ws.Range("E70").Value2 = "ROLES"
$ws.Range("F70").Value2 = "ROLES"
$ws.Range("G70").Value2 = 0

# --- Update the active selection to match the committed view state ---
$ws.Range("F63").Select()
